$d = $word.ActiveDocument

# 1. Replace the ID placeholder text in the first paragraph
$d.Content.Find.Execute("**ID__AFFARS_mp_5315_3_topic_3__ID**", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_MP_5315_3_1_1__ID**", 2)

# 2. Remove the now-orphaned trailing-space run that followed the ID text
#    (the paragraph's last character is the paragraph mark; the one before
#    it is the stray space run left over from the old two-run layout)
$p1 = $d.Paragraphs.Item(1)
$charCount = $p1.Range.Characters.Count
$spaceChar = $p1.Range.Characters.Item($charCount - 1)
$spaceChar.Delete()

# 3. Update paragraph formatting: indent and paragraph border (space-only, no visible lines)
$pPr = $p1.Range.ParagraphFormat
$pPr.LeftIndent = 11.25

$borders = $pPr.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
